# Updating testcases in master:
# Mark all the remaining "NO" testcases in column E as "Yes" (the row that
# was already "Yes" - E29 - is left untouched, matching the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

$ws.Range("E3:E28").Value = "Yes"
$ws.Range("E30:E39").Value = "Yes"

# Reflect the final cursor position left in the sheet after the edit.
$ws.Range("E4").Select()
